# Actualización automática 2025-09-22 14:30:09
$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual    = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento    = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO ---
# M2: PORCELANATO for ALTAMIRANO MARCATOMA EDISON PAULINO
$wsVentasPorGrupo.Range("M2").Value = 444.52
# P20: NO RESURTIBLES for MANCHENO PINO HERVIN SANTIAGO
$wsVentasPorGrupo.Range("P20").Value = 359.95

# --- VENTA MENSUAL ---
# F2: septiembre for ALTAMIRANO MARCATOMA EDISON PAULINO
$wsVentaMensual.Range("F2").Value = 444.52
# F20: septiembre for MANCHENO PINO HERVIN SANTIAGO
$wsVentaMensual.Range("F20").Value = 2703.86
# F35: septiembre total
$wsVentaMensual.Range("F35").Value = 16792.83

# --- CUMPLIMIENTO MENSUAL ---
# Row 8: NO RESURTIBLES
$wsCumplimiento.Range("D8").Value = 359.95
$wsCumplimiento.Range("E8").Value = 120.267743214072
$wsCumplimiento.Range("F8").Value = 0.7495558110595282

# Row 12: PORCELANATO
$wsCumplimiento.Range("D12").Value = 10240.55
$wsCumplimiento.Range("E12").Value = 12193.2053751766
$wsCumplimiento.Range("F12").Value = 0.4564795251057865

# Row 15: TOTAL
$wsCumplimiento.Range("D15").Value = 17050.06
$wsCumplimiento.Range("E15").Value = 21692.95881339593
$wsCumplimiento.Range("F15").Value = 0.4400808331978691

# Column D width widened slightly (matches observed autofit after value-length change)
# Note: stored OOXML <col width> = ColumnWidth + 5/6 for this engine's font metrics,
# so we request 14 - 5/6 to land on a stored width of exactly 14.
$wsCumplimiento.Columns.Item(4).ColumnWidth = 13.166666666666666
